$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values (Price / Volume columns) need a temporary
# Text number-format so Excel stores them as literal strings instead of
# coercing "258.58" -> 258.58 (number) or "5.48%" -> 0.0548 (percentage).
# Resetting the Style back to "Normal" afterwards drops the temporary
# format so the cell keeps its original (default) style index.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "258.58"
Set-TextValue $ws.Range("E2") "5.48%"
Set-TextValue $ws.Range("D3") "28.06"
Set-TextValue $ws.Range("E3") "-3.85%"
Set-TextValue $ws.Range("D4") "5.217"
Set-TextValue $ws.Range("E4") "-0.83%"
Set-TextValue $ws.Range("D5") "0.05927"
Set-TextValue $ws.Range("E5") "3.88%"
Set-TextValue $ws.Range("D6") "6.707"
Set-TextValue $ws.Range("E6") "1.39%"
Set-TextValue $ws.Range("D7") "0.8745"
Set-TextValue $ws.Range("E7") "2.72%"
Set-TextValue $ws.Range("D8") "0.9893"
Set-TextValue $ws.Range("E8") "15.33%"
Set-TextValue $ws.Range("D9") "0.1416"
Set-TextValue $ws.Range("E9") "3.40%"
Set-TextValue $ws.Range("D10") "0.07193"
Set-TextValue $ws.Range("E10") "2.05%"
Set-TextValue $ws.Range("D11") "0.03200"
Set-TextValue $ws.Range("E11") "0.27%"
Set-TextValue $ws.Range("D12") "0.09240"
Set-TextValue $ws.Range("E13") "1.18%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D14") "0.01065"
Set-TextValue $ws.Range("E14") "1,680.88%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.006032"
Set-TextValue $ws.Range("E15") "1.75%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.496"
Set-TextValue $ws.Range("E16") "0.07%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "3.228"
Set-TextValue $ws.Range("E17") "1.13%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D18") "2.205"
Set-TextValue $ws.Range("E18") "1.36%"
Set-TextValue $ws.Range("D19") "0.3122"
Set-TextValue $ws.Range("E19") "-1.24%"
Set-TextValue $ws.Range("D20") "0.03714"
Set-TextValue $ws.Range("E20") "12.01%"
Set-TextValue $ws.Range("E21") "1.60%"
Set-TextValue $ws.Range("D22") "3.519"
Set-TextValue $ws.Range("E22") "0.29%"
Set-TextValue $ws.Range("D23") "0.04201"
Set-TextValue $ws.Range("E23") "2.43%"
Set-TextValue $ws.Range("D24") "0.1397"
Set-TextValue $ws.Range("E24") "1.28%"
Set-TextValue $ws.Range("D25") "0.001216"
Set-TextValue $ws.Range("E25") "-0.73%"
Set-TextValue $ws.Range("D26") "0.004554"
Set-TextValue $ws.Range("E26") "9.92%"
Set-TextValue $ws.Range("E27") "-0.12%"
Set-TextValue $ws.Range("D28") "0.0001936"
Set-TextValue $ws.Range("E28") "33.57%"
Set-TextValue $ws.Range("D40") "0.03843"
Set-TextValue $ws.Range("E40") "2.41%"
Set-TextValue $ws.Range("D41") "0.005456"
Set-TextValue $ws.Range("E41") "5.39%"
Set-TextValue $ws.Range("D42") "0.1104"
Set-TextValue $ws.Range("E42") "3.94%"
Set-TextValue $ws.Range("D43") "0.002298"
Set-TextValue $ws.Range("E43") "-6.19%"
Set-TextValue $ws.Range("D44") "0.01073"
Set-TextValue $ws.Range("E44") "14.75%"
Set-TextValue $ws.Range("E45") "2.75%"
Set-TextValue $ws.Range("D46") "0.00000000749"
Set-TextValue $ws.Range("E46") "-0.14%"
Set-TextValue $ws.Range("D47") "0.08542"
Set-TextValue $ws.Range("E47") "13.84%"
Set-TextValue $ws.Range("D48") "0.002136"
Set-TextValue $ws.Range("E48") "-12.52%"
Set-TextValue $ws.Range("D49") "0.00002098"
Set-TextValue $ws.Range("E49") "-0.14%"
Set-TextValue $ws.Range("D50") "0.0001998"
Set-TextValue $ws.Range("E50") "-0.14%"
